# TournRPG-183 バフスキルの実装
# Rename a handful of skill names/details on the "skill" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill")

# SKILL001: 大斬り -> かみつき
$ws.Range("B3").Value = "かみつき"

# SKILL005: HP回復 -> ヒーリング
$ws.Range("B7").Value = "ヒーリング"

# SKILL006: 攻撃アップ -> パワーアップ, 攻撃力アップ -> 攻撃力をアップする
$ws.Range("B8").Value = "パワーアップ"
$ws.Range("T8").Value = "攻撃力をアップする"
